$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("4A", 0.853),
    @("4P", 0.893),
    @("4S", 0.9508333329999999),
    @("5P", 0.99025),
    @("5W", 0.873),
    @("8F", 0.8784999999999999),
    @("8H", 0.9105),
    @("8K", 0.9385),
    @("10H", 0.9025),
    @("20H", 0.9695),
    @("55I", 0.916),
    @("71K", 0.9895)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
